$d = $word.ActiveDocument

$replacements = @(
    @("66÷3=22, 0", "12÷4=3, 0"),
    @("70÷4=17, 2", "37÷8=4, 5"),
    @("68÷7=9, 5", "93÷8=11, 5"),
    @("78÷8=9, 6", "82÷9=9, 1"),
    @("67÷4=16, 3", "25÷9=2, 7"),
    @("32÷6=5, 2", "73÷2=36, 1"),
    @("19÷6=3, 1", "27÷8=3, 3"),
    @("95÷2=47, 1", "51÷9=5, 6"),
    @("84÷6=14, 0", "36÷4=9, 0"),
    @("90÷4=22, 2", "23÷5=4, 3"),
    @("78÷3=26, 0", "88÷6=14, 4"),
    @("59÷9=6, 5", "69÷6=11, 3"),
    @("66÷9=7, 3", "28÷7=4, 0"),
    @("25÷6=4, 1", "64÷3=21, 1"),
    @("86÷8=10, 6", "33÷5=6, 3"),
    @("64÷5=12, 4", "70÷2=35, 0"),
    @("53÷6=8, 5", "18÷6=3, 0"),
    @("77÷3=25, 2", "55÷5=11, 0"),
    @("98÷4=24, 2", "45÷3=15, 0"),
    @("10÷3=3, 1", "50÷6=8, 2"),
    @("72÷7=10, 2", "33÷4=8, 1"),
    @("30÷6=5, 0", "19÷9=2, 1"),
    @("55÷9=6, 1", "23÷3=7, 2"),
    @("53÷2=26, 1", "14÷6=2, 2"),
    @("51÷8=6, 3", "77÷5=15, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
